$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 27, pushing existing rows 27..51 down to 29..53.
$ws.Rows.Item(27).Insert()
$ws.Rows.Item(27).Insert()

# New row 27: Perejil, Primera, week of 2023-05-03 (serial 45049)
$ws.Range("A27").Value = 7
$ws.Range("B27").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C27").Value = "Ñuble"
$ws.Range("D27").Value = 45049
$ws.Range("E27").Value = 16
$ws.Range("F27").Value = 100112044
$ws.Range("G27").Value = "Perejil"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 100
$ws.Range("K27").Value = 1200
$ws.Range("L27").Value = 1200
$ws.Range("M27").Value = 1200
$ws.Range("N27").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O27").Value = "Región del Maule"
$ws.Range("P27").Value = 1200
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = "Hortaliza"

# New row 28: Perejil, Segunda, week of 2023-05-03 (serial 45049)
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 45049
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 100112044
$ws.Range("G28").Value = "Perejil"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Segunda"
$ws.Range("J28").Value = 150
$ws.Range("K28").Value = 1000
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 1000
$ws.Range("N28").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O28").Value = "Región del Maule"
$ws.Range("P28").Value = 1000
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = "Hortaliza"
